$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 22:05"

# --- Update Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1786355
$ws.Range("C4").Value = 17894
$ws.Range("D4").Value = 502305
$ws.Range("E4").Value = 1179915
$ws.Range("G4").Value = 805
$ws.Range("H4").Value = 104135

# --- Update Brasil (row 5) ---
$ws.Range("B5").Value = 450079
$ws.Range("C5").Value = 11267
$ws.Range("E5").Value = 229622
$ws.Range("G5").Value = 512
$ws.Range("H5").Value = 27276

# --- Update Alemania (row 11) ---
$ws.Range("B11").Value = 183019
$ws.Range("C11").Value = 567
$ws.Range("E11").Value = 10325
$ws.Range("G11").Value = 24
$ws.Range("H11").Value = 8594

# --- Update India (row 12) ---
$ws.Range("B12").Value = 173491
$ws.Range("C12").Value = 8105
$ws.Range("D12").Value = 82627
$ws.Range("E12").Value = 85884

# --- Update Sudafrica (row 33) ---
$ws.Range("B33").Value = 29240
$ws.Range("C33").Value = 1837
$ws.Range("D33").Value = 15093
$ws.Range("E33").Value = 13536
$ws.Range("G33").Value = 34
$ws.Range("H33").Value = 611

# --- Update Camerun (row 69) ---
$ws.Range("D69").Value = 3326
$ws.Range("E69").Value = 1933
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 177

# --- Swap "Santa Lucia" / "Belice" ordering ---
# Row 200 was Santa Lucia, row 201 was Belice; the country list order
# swaps so Belice now sorts before Santa Lucia (same totals columns
# follow the name to its new row).
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# --- Swap "San Bartolome" / "Bonaire, San Eustaquio y Saba" ordering ---
# Row 215 was San Bartolome, row 216 was Bonaire, San Eustaquio y Saba.
# Both rows carry identical totals (6,0,6,0,0,0,0), so only the names
# need to trade places -- no numeric columns change.
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"
